# Updated capital structure database
#
# The comdirect bank AG (XTRA:COM) row is removed from the Germany / Bank (Money
# Center) peer set. The remaining companies keep their relative order except that
# Commerzbank AG and HSBC Trinkaus & Burkhardt AG swap rows, and every recalculated
# metric column (growth, margins, cash, debt, coverage ratios, etc.) is refreshed
# to reflect the updated peer-group composition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the comdirect bank AG (XTRA:COM) row entirely (row 7).
# Rows 2-6 keep their row numbers; only row 7 disappears.
$ws.Rows(7).Delete()


# Row 2: industry aggregate / comparison row (company count 5 -> 4 after removal).
$ws.Range("B2").Value = "4"
$ws.Range("D2").Value = -0.00515
$ws.Range("E2").Value = -0.0143
$ws.Range("F2").Value = -0.06860000000000001
$ws.Range("I2").Value = [double]"-6.528183262934564e-05"
$ws.Range("J2").Value = [double]"-4.451589583863478e-05"
$ws.Range("K2").Value = -7.59999999999998
$ws.Range("L2").Value = -0.0007883572088005538
$ws.Range("M2").Value = 46.03279999999999
$ws.Range("N2").Value = 0.003940658305868253
$ws.Range("O2").Value = -6.056947368421068
$ws.Range("P2").Value = 46.03279999999999
$ws.Range("Q2").Value = 0.003940658305868253
$ws.Range("R2").Value = -6.056947368421068
$ws.Range("U2").Value = 110550.2
$ws.Range("V2").Value = 9.463699011257116
$ws.Range("W2").Value = 0.05676154332458941
$ws.Range("X2").Value = 0.07960543238549264
$ws.Range("Y2").Value = -0.02284388906090323
$ws.Range("Z2").Value = 0.1164440072923984
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.03042416642154142
$ws.Range("AC2").Value = -0.03180571047748443
$ws.Range("AD2").Value = 117654.3
$ws.Range("AE2").Value = 8.346682255483405
$ws.Range("AF2").Value = 117662.6466822555
$ws.Range("AG2").Value = 7112.446682255468
$ws.Range("AH2").Value = 0.9096866746610763
$ws.Range("AI2").Value = 0.7390923309708268
$ws.Range("AJ2").Value = 0.3784434851552905
$ws.Range("AK2").Value = 0.1462002988224896
$ws.Range("AN2").Value = 113129.1346153846
$ws.Range("AP2").Value = 6838.891040630258

# Row 3: ProCredit Holding AG & Co. KGaA (XTRA:PCZ) - refreshed metrics.
$ws.Range("D3").Value = -0.075
$ws.Range("E3").Value = -0.07719999999999999
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 51.1
$ws.Range("L3").Value = 0.1834829443447038
$ws.Range("M3").Value = 20.7328
$ws.Range("N3").Value = 0.03862295081967213
$ws.Range("O3").Value = 0.4057299412915851
$ws.Range("P3").Value = 20.7328
$ws.Range("Q3").Value = 0.03862295081967213
$ws.Range("R3").Value = 0.4057299412915851
$ws.Range("U3").Value = 145.1
$ws.Range("V3").Value = 0.2703055141579732
$ws.Range("W3").Value = 0.06001174398120963
$ws.Range("X3").Value = 0.08690702611487093
$ws.Range("Y3").Value = -0.0268952821336613
$ws.Range("Z3").Value = 0.1176893171061528
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.02986801783656699
$ws.Range("AC3").Value = -0.02986801783656699
$ws.Range("AD3").Value = 1863.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1863.8
$ws.Range("AG3").Value = 1718.7
$ws.Range("AH3").Value = 0.7763892360243273
$ws.Range("AI3").Value = 0.6690118094691123
$ws.Range("AJ3").Value = 0.7620039902460652
$ws.Range("AK3").Value = 0.6508255074219933
$ws.Cells.Item(3,40).ClearContents()  # AN3
$ws.Cells.Item(3,42).ClearContents()  # AP3

# Row 4: now HSBC Trinkaus & Burkhardt AG (DB:TUB) (moved up from row 5).
$ws.Range("B4").Value = "HSBC Trinkaus & Burkhardt AG (DB:TUB)"
$ws.Range("D4").Value = 0.0268
$ws.Range("E4").Value = -0.0143
$ws.Cells.Item(4,6).ClearContents()  # F4
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 147.9
$ws.Range("L4").Value = 0.1649565023421816
$ws.Range("M4").Value = 25.3
$ws.Range("N4").Value = 0.008485946199771919
$ws.Range("O4").Value = 0.1710615280594997
$ws.Range("P4").Value = 25.3
$ws.Range("Q4").Value = 0.008485946199771919
$ws.Range("R4").Value = 0.1710615280594997
$ws.Range("U4").Value = 11827
$ws.Range("V4").Value = 3.966928288723418
$ws.Range("W4").Value = 0.05351134266796918
$ws.Range("X4").Value = 0.04312265916664064
$ws.Range("Y4").Value = 0.01038868350132854
$ws.Range("Z4").Value = -0.2371079494367165
$ws.Range("AA4").Value = -0
$ws.Range("AB4").Value = 0.0308265655591313
$ws.Range("AC4").Value = -0.0308265655591313
$ws.Range("AD4").Value = 2108.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 2108.5
$ws.Range("AG4").Value = -9718.5
$ws.Range("AH4").Value = 0.4142517534725633
$ws.Range("AI4").Value = 0.4265196722969556
$ws.Range("AJ4").Value = 1.442534621721512
$ws.Range("AK4").Value = 1.411854434517324
$ws.Cells.Item(4,40).ClearContents()  # AN4
$ws.Cells.Item(4,42).ClearContents()  # AP4

# Row 5: now Commerzbank AG (XTRA:CBK) (moved down from row 4).
$ws.Range("B5").Value = "Commerzbank AG (XTRA:CBK)"
$ws.Range("D5").Value = -0.0371
$ws.Cells.Item(5,5).ClearContents()  # E5
$ws.Range("F5").Value = -0.06860000000000001
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -233.7
$ws.Range("L5").Value = -0.02784264201296225
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = 0
$ws.Cells.Item(5,20).ClearContents()  # T5
$ws.Range("U5").Value = 98258.60000000001
$ws.Range("V5").Value = 12.18001289170964
$ws.Range("W5").Value = -0.007284799177070182
$ws.Range("X5").Value = 0.2546193460291942
$ws.Range("Y5").Value = -0.2619041452062644
$ws.Range("Z5").Value = 0.09976762568122523
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.03278485539583757
$ws.Range("AC5").Value = -0.03278485539583757
$ws.Range("AD5").Value = 113445.3
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 113445.3
$ws.Range("AG5").Value = 15186.7
$ws.Range("AH5").Value = 0.9336101224153894
$ws.Range("AI5").Value = 0.7513028309627737
$ws.Range("AJ5").Value = 0.6530818486361427
$ws.Range("AK5").Value = 0.2879568444903725
$ws.Cells.Item(5,40).ClearContents()  # AN5
$ws.Cells.Item(5,42).ClearContents()  # AP5

# Row 6: Merkur Bank KGaA (XTRA:MBK) - refreshed metrics.
$ws.Range("D6").Value = 0.145
$ws.Range("E6").Value = 0.451
$ws.Range("I6").Value = -0.00878961523878046
$ws.Range("J6").Value = -0.006834683573603425
$ws.Range("K6").Value = 27.1
$ws.Range("L6").Value = 0.3784916201117319
$ws.Range("U6").Value = 319.5
$ws.Range("V6").Value = 3.324661810613944
$ws.Range("W6").Value = 0.1966618287373004
$ws.Range("X6").Value = 0.07230383865611435
$ws.Range("Y6").Value = 0.1243579900811861
$ws.Range("Z6").Value = 0.9855921533786987
$ws.Range("AA6").Value = -0.006736210500969819
$ws.Range("AB6").Value = 0.03002176728395155
$ws.Range("AC6").Value = -0.03675797778492137
$ws.Range("AD6").Value = 236.7
$ws.Range("AE6").Value = 8.346682255483405
$ws.Range("AF6").Value = 245.0466822554834
$ws.Range("AG6").Value = -74.4533177445166
$ws.Range("AH6").Value = 0.7183029910634421
$ws.Range("AI6").Value = 0.5198862991522258
$ws.Range("AJ6").Value = -3.439479402237567
$ws.Range("AK6").Value = -0.4903190286321055
$ws.Range("AN6").Value = 227.5961538461538
$ws.Range("AP6").Value = -71.58972860049673
